# Corrige os dados e atualiza o cabecalho da planilha PNAD 2009 (roubo/furto)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A segunda linha de cabecalho (B2:F2) tinha rotulos "unnamed: 1_level_1" e
# "unnamed: 5_level_1" que eram, na verdade, repeticoes do rotulo "total".
$ws.Cells.Item(2, 2).Value2 = "total"
$ws.Cells.Item(2, 6).Value2 = "total"

# As linhas 5 ("situacao do domicilio") e 8 ("grandes regioes e unidades da
# federacao") eram apenas cabecalhos de secao sem nenhum dado associado.
# Remove-las para que os dados fiquem corretamente alinhados com seus
# rotulos (ex.: "urbana" e "rural" passam a ter os valores numericos na
# mesma linha). Exclui de baixo para cima para nao desalinhar os indices.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
